$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 1003, shifting existing rows 1003:1079 down to 1004:1080
$ws.Rows.Item(1003).Insert()

# Populate the newly inserted row with this week's data point
$ws.Range("A1003").Value = 6
$ws.Range("B1003").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1003").Value = "Metropolitana"
$ws.Range("D1003").Value = 45021
$ws.Range("E1003").Value = 13
$ws.Range("F1003").Value = 100112003
$ws.Range("G1003").Value = "Ajo"
$ws.Range("H1003").Value = "Chino"
$ws.Range("I1003").Value = "Primera"
$ws.Range("J1003").Value = 900
$ws.Range("K1003").Value = 13000
$ws.Range("L1003").Value = 14000
$ws.Range("M1003").Value = 13333
$ws.Range("N1003").Value = "$/caja 10 kilos"
$ws.Range("O1003").Value = "China"
$ws.Range("P1003").Value = 1333
$ws.Range("Q1003").Value = 10
$ws.Range("R1003").Value = "Hortaliza"

# Make sure the date cell keeps the same number format as the rest of column D
$ws.Range("D1003").NumberFormat = $ws.Range("D1004").NumberFormat
